$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.247.07'
$ws.Range("E2").Value = '  +1.53%  '
$ws.Range("D3").Value = '2.349.59'
$ws.Range("E3").Value = '  +2.34%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '303.26'
$ws.Range("E5").Value = '  +0.81%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '95.62'
$ws.Range("E6").Value = '  -0.23%  '
$ws.Range("E7").Value = '  -0.20%  '
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.498'
$ws.Range("E9").Value = '  +0.72%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '34.28'
$ws.Range("E10").Value = '  -1.11%  '
$ws.Range("E11").Value = '  +0.37%  '
$ws.Range("B12").Value = 'Chainlink'
$ws.Range("C12").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '18.55'
$ws.Range("E12").Value = '  -3.39%  '
$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.121'
$ws.Range("E13").Value = '  +2.39%  '
$ws.Range("E14").Value = '  -0.04%  '
$ws.Range("D15").Value = '2.719.11'
$ws.Range("E15").Value = '  +2.55%  '
$ws.Range("D16").Value = '2.357.94'
$ws.Range("E16").Value = '  +2.83%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.798'
$ws.Range("E17").Value = '  +2.05%  '
$ws.Range("D18").Value = '43.220.75'
$ws.Range("E18").Value = '  +1.54%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.25'
$ws.Range("E19").Value = '  +0.48%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.21'
$ws.Range("E20").Value = '  +3.44%  '
$ws.Range("E21").Value = '  +0.37%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '68.25'
$ws.Range("E22").Value = '  +1.04%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '236.04'
$ws.Range("E23").Value = '  +0.50%  '
$ws.Range("E24").Value = '  -1.04%  '
$ws.Range("E25").Value = '  +0.04%  '
$ws.Range("E26").Value = '  +0.88%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '24.67'
$ws.Range("E28").Value = '  +14.89%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.17'
$ws.Range("E29").Value = '  +1.65%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '31.45'
$ws.Range("E30").Value = '  -2.37%  '
$ws.Range("E31").Value = '  +0.00%  '
$ws.Range("E32").Value = '  +0.97%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0727'
$ws.Range("E33").Value = '  +4.29%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '17.31'
$ws.Range("E34").Value = '  -0.54%  '
$ws.Range("E35").Value = '  +5.50%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.39'
$ws.Range("E36").Value = '  -0.59%  '
$ws.Range("B38").Value = 'EnergySwap'
$ws.Range("C38").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '22.78'
$ws.Range("E38").Value = '  +16.32%  '
$ws.Range("B39").Value = 'Kaspa'
$ws.Range("C39").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.101'
$ws.Range("E39").Value = '  +0.92%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.75'
$ws.Range("E40").Value = '  +1.28%  '
$ws.Range("E41").Value = '  -0.04%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '113.55'
$ws.Range("E42").Value = '  -30.82%  '
$ws.Range("D43").Value = '1.942.48'
$ws.Range("E43").Value = '  -0.62%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0282'
$ws.Range("E44").Value = '  +0.99%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.01'
$ws.Range("E45").Value = '  -4.67%  '
$ws.Range("E46").Value = '  +2.51%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.73'
$ws.Range("E47").Value = '  -0.59%  '
$ws.Range("D48").Value = '2.582.24'
$ws.Range("E48").Value = '  +2.44%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '53.09'
$ws.Range("E49").Value = '  +0.01%  '
$ws.Range("E50").Value = '  -3.63%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '72.24'
$ws.Range("E51").Value = '  +0.96%  '
